# "se agrega procolo y pines al informe"
# Adds new analog-pin assignments (A3, A4) for the two temperature actuators,
# and a new row for "Coolers 1 y 2" using pin A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: "Actuador temperatura calentar" now uses analog pin A3 instead
#     of digital pin 0 ---
$ws.Range("B10").ClearContents()
$ws.Range("C10").Value = "A3"
$ws.Range("C4").Copy()
$ws.Range("C10").PasteSpecial(-4122)   # xlPasteFormats

# --- Row 11: "Actuador temperatura enfriar" now uses analog pin A4 instead
#     of digital pin 1 ---
$ws.Range("B11").ClearContents()
$ws.Range("C11").Value = "A4"
$ws.Range("C4").Copy()
$ws.Range("C11").PasteSpecial(-4122)   # xlPasteFormats

# --- New row 14: "Coolers 1 y 2" using analog pin A2 ---
$ws.Range("A13:D13").Copy()
$ws.Range("A14:D14").PasteSpecial(-4122)   # xlPasteFormats, mirrors row 13 look
$ws.Range("A14").Value = "Coolers 1 y 2"
$ws.Range("B14").ClearContents()
$ws.Range("C14").Value = "A2"
$ws.Range("C4").Copy()
$ws.Range("C14").PasteSpecial(-4122)   # xlPasteFormats (right-aligned pin style)
$ws.Range("D14").ClearContents()

$excel.CutCopyMode = 0

# --- Update the selection shown when the workbook is next opened ---
$ws.Range("A16").Select()
